# features selection, 5 runs
# Update AdaBoostClassifier_LogisticRegression results sheet: row 2 gets new
# metrics, rows 3-5 are refreshed/reshuffled experiment rows, and row 6 is a
# brand-new run (feature-selector experiments). Dimension grows to A1:K6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @'
Pipeline(steps=[('scaler', StandardScaler()),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(C=3,
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver='saga'),
                                    random_state=42))])
'@
$ws.Range("B2").Value = 0.7047619047619047
$ws.Range("C2").Value = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': StandardScaler(), 'model__n_estimators': 50, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': None, 'model__estimator__C': 3}
'@
$ws.Range("D2").Value = 0.3333333333333333
$ws.Range("E2").Value = @'
[1 1 0 0 1 0 0 0 0 1 0 1]
'@
$ws.Range("F2").Value = @'
[0 0 1 0 0 1 1 0 1 1 1 1]
'@
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.6176666666666667
$ws.Range("I2").Value = 0.02250861618237601
$ws.Range("J2").Value = 0.5615238095238095
$ws.Range("K2").Value = 0.0481743073251333

$ws.Range("A3").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(C=3,
                                                                 class_weight='balanced',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver='saga'),
                                    random_state=42))])
'@
$ws.Range("B3").Value = 0.6285714285714284
$ws.Range("C3").Value = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 3}
'@
$ws.Range("D3").Value = 0.8571428571428571
$ws.Range("E3").Value = @'
[1 1 0 1 0 0 1 0 1 1 1 0]
'@
$ws.Range("F3").Value = @'
[0 1 0 1 1 0 1 0 1 1 1 0]
'@
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.5593112244897959
$ws.Range("I3").Value = 0.02097903554945518
$ws.Range("J3").Value = 0.5183673469387754
$ws.Range("K3").Value = 0.04995927608673544

$ws.Range("A4").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight='balanced',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver='saga'),
                                    random_state=42))])
'@
$ws.Range("B4").Value = 0.6285714285714287
$ws.Range("C4").Value = @'
{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.0001}
'@
$ws.Range("D4").Value = 0.5714285714285715
$ws.Range("E4").Value = @'
[1 0 1 1 1 1 0 1 0 1 0 1]
'@
$ws.Range("F4").Value = @'
[0 0 0 1 1 1 1 0 1 0 0 1]
'@
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6042424242424242
$ws.Range("I4").Value = 0.02880239252578566
$ws.Range("J4").Value = 0.5478787878787879
$ws.Range("K4").Value = 0.06461422295917255

$ws.Range("A5").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1',
                                                     random_state=42))),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight='balanced',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver='saga'),
                                    random_state=42))])
'@
$ws.Range("B5").Value = 0.6095238095238095
$ws.Range("C5").Value = @'
{'selector': SelectFromModel(estimator=LinearSVC(dual=False, penalty='l1', random_state=42)), 'scaler': None, 'model__n_estimators': 50, 'model__estimator__solver': 'saga', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.0001}
'@
$ws.Range("D5").Value = 0.75
$ws.Range("E5").Value = @'
[1 1 0 0 0 0 1 0 1 1 1 1]
'@
$ws.Range("F5").Value = @'
[1 1 1 0 0 1 1 1 1 0 1 1]
'@
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.5905643738977071
$ws.Range("I5").Value = 0.02121079951592237
$ws.Range("J5").Value = 0.51657848324515
$ws.Range("K5").Value = 0.05448520659315167

$ws.Range("A6").Value = @'
Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45a99f0f40>),
                ('model',
                 AdaBoostClassifier(estimator=LogisticRegression(C=0.0001,
                                                                 class_weight='balanced',
                                                                 max_iter=1000,
                                                                 random_state=42,
                                                                 solver='liblinear'),
                                    n_estimators=5, random_state=42))])
'@
$ws.Range("B6").Value = 0.638095238095238
$ws.Range("C6").Value = @'
{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c418fd0>, 'scaler': None, 'model__n_estimators': 5, 'model__estimator__solver': 'liblinear', 'model__estimator__penalty': 'l2', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 0.0001}
'@
$ws.Range("D6").Value = 0.6666666666666666
$ws.Range("E6").Value = @'
[1 1 1 1 0 0 0 0 1 1 0 0]
'@
$ws.Range("F6").Value = @'
[1 1 1 0 0 1 0 1 1 0 0 0]
'@
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.5943027210884353
$ws.Range("I6").Value = 0.02016137469583393
$ws.Range("J6").Value = 0.5346938775510204
$ws.Range("K6").Value = 0.05071880798722873
